# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Rafflesia_Profits workbook (per-sheet Leve profit calculation columns H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1571
$ws.Range("I19").Value = 1749.5
$ws.Range("J19").Value = 1333
$ws.Range("K19").Value = 1749.5
$ws.Range("L19").Value = 1333
$ws.Range("M19").Value = -1574.5
$ws.Range("N19").Value = -1683
$ws.Range("H53").Value = 1787
$ws.Range("I53").Value = 1549.3334
$ws.Range("J53").Value = 2500
$ws.Range("K53").Value = 1549.3334
$ws.Range("L53").Value = 2500
$ws.Range("M53").Value = -912.3334
$ws.Range("N53").Value = -3774
$ws.Range("H98").Value = 621.2
$ws.Range("I98").Value = 579.1111
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 579.1111
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 918.8889
$ws.Range("N98").Value = -3996
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("M116").Value = 3000
$ws.Range("N116").Value = -9884
$ws.Range("H122").Value = 621.2
$ws.Range("I122").Value = 579.1111
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 1737.3333
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 712.6667000000002
$ws.Range("N122").Value = -7900
$ws.Range("H138").Value = 2576.4707
$ws.Range("I138").Value = 733.3333
$ws.Range("J138").Value = 2971.4285
$ws.Range("K138").Value = 2199.9999
$ws.Range("L138").Value = 8914.2855
$ws.Range("M138").Value = 2940.0001
$ws.Range("N138").Value = -19194.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1800
$ws.Range("I20").Value = 1800
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1800
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1553
$ws.Range("H82").Value = 5197.143
$ws.Range("I82").Value = 5197.143
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 5197.143
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -4814.143
$ws.Range("H85").Value = 5197.143
$ws.Range("I85").Value = 5197.143
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 5197.143
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -3871.143
$ws.Range("H99").Value = 1152.5
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1152.5
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1152.5
$ws.Range("N99").Value = -4148.5
$ws.Range("H135").Value = 96000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 96000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 96000
$ws.Range("N135").Value = -106140
$ws.Range("H139").Value = 80000
$ws.Range("I139").Value = 80000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 80000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -74860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 20.8
$ws.Range("I7").Value = 9.142858
$ws.Range("J7").Value = 48
$ws.Range("K7").Value = 9.142858
$ws.Range("L7").Value = 48
$ws.Range("M7").Value = 103.857142
$ws.Range("N7").Value = -274
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H69").Value = 4875
$ws.Range("I69").Value = 4875
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 4875
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -4126
$ws.Range("H72").Value = 4875
$ws.Range("I72").Value = 4875
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 14625
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -10881
$ws.Range("H132").Value = 10412.4
$ws.Range("I132").Value = 10412.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 31237.2
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -28707.2
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1001
$ws.Range("I80").Value = 1001
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3003
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2067
$ws.Range("H83").Value = 1001
$ws.Range("I83").Value = 1001
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9009
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -4329
$ws.Range("H113").Value = 6013.1665
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 7035.8
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 21107.4
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -25447.4
$ws.Range("H122").Value = 761.8
$ws.Range("I122").Value = 749.5
$ws.Range("J122").Value = 770
$ws.Range("K122").Value = 6745.5
$ws.Range("L122").Value = 6930
$ws.Range("M122").Value = -4295.5
$ws.Range("N122").Value = -11830
$ws.Range("H123").Value = 1500
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 2000
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 6000
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -10900
$ws.Range("H132").Value = 1081.6666
$ws.Range("I132").Value = 1150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 40.76923
$ws.Range("I2").Value = 59.875
$ws.Range("J2").Value = 10.2
$ws.Range("K2").Value = 59.875
$ws.Range("L2").Value = 10.2
$ws.Range("M2").Value = 53.125
$ws.Range("N2").Value = -236.2
$ws.Range("H43").Value = 3344.6667
$ws.Range("I43").Value = 3344.6667
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3344.6667
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3193.6667
$ws.Range("H47").Value = 15000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 15000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16136
$ws.Range("H63").Value = 29996.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 29996.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 29996.5
$ws.Range("N63").Value = -31368.5
$ws.Range("H66").Value = 29996.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 29996.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 89989.5
$ws.Range("N66").Value = -96853.5
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H97").Value = 533.3333
$ws.Range("I97").Value = 490
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 490
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = 6
$ws.Range("N97").Value = -1742
$ws.Range("H122").Value = 817
$ws.Range("I122").Value = 817
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2451
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1888
$ws.Range("H16").Value = 15666.667
$ws.Range("I16").Value = 15666.667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 15666.667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -15496.667
$ws.Range("H22").Value = 1116.6666
$ws.Range("I22").Value = 925
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 925
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -630
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1116.6666
$ws.Range("I27").Value = 925
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 925
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -818
$ws.Range("N27").Value = -1714
$ws.Range("H40").Value = 4199.8
$ws.Range("I40").Value = 3999.75
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3999.75
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3863.75
$ws.Range("N40").Value = -5272
$ws.Range("H46").Value = 775
$ws.Range("I46").Value = 775
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 775
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = -587
$ws.Range("N46").Value = 0
$ws.Range("H93").Value = 1100
$ws.Range("I93").Value = 1100
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1100
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 148
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3530
$ws.Range("H136").Value = 57197
$ws.Range("I136").Value = 24001.273
$ws.Range("J136").Value = 97769.55499999999
$ws.Range("K136").Value = 72003.819
$ws.Range("L136").Value = 293308.665
$ws.Range("M136").Value = -69453.819
$ws.Range("N136").Value = -298408.665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = 0
